# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.744.28'
$ws.Range("E2").Value = '  +2.00%  '
$ws.Range("D3").Value = '1.574.79'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''213.44'
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").Value = '''0.491'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''44.88'
$ws.Range("E8").Value = '  +2.20%  '
$ws.Range("D9").Value = '''24.16'
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").Value = '''0.0891'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '1.801.29'
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("D14").Value = '1.582.40'
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D16").Value = '28.754.09'
$ws.Range("E16").Value = '  +1.84%  '
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '''62.51'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").Value = '''231.47'
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("D20").Value = '''7.40'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").Value = '0.0₃0693'
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  -4.71%  '
$ws.Range("D24").Value = '''9.18'
$ws.Range("E24").Value = '  -1.26%  '
$ws.Range("D25").Value = '''2.10'
$ws.Range("E25").Value = '  +7.88%  '
$ws.Range("D26").Value = '''152.27'
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").Value = '''15.02'
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("E28").Value = '  -1.37%  '
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = '''0.0483'
$ws.Range("E31").Value = '  +2.70%  '
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").Value = '''3.21'
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("D35").Value = '1.397.64'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D37").Value = '''1.54'
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("E39").Value = '  +2.93%  '
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").Value = '''0.526'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D45").Value = '''0.0472'
$ws.Range("E45").Value = '  +3.30%  '
$ws.Range("D46").Value = '''5.51'
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("D47").Value = '''0.961'
$ws.Range("E47").Value = '  -2.01%  '
$ws.Range("D48").Value = '''63.29'
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("D49").Value = '1.713.18'
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").Value = '''86.70'
$ws.Range("D51").Value = '0.0₆0103'
$ws.Range("E51").Value = '  +0.57%  '

# Rows 43 and 44 swap positions (ARBITRUM <-> RenderToken) with refreshed values
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '''0.794'
$ws.Range("E43").Value = '  -2.04%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''1.90'
$ws.Range("E44").Value = '  +1.08%  '
